# Update Name of Algo
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B7").Value  = 5.697000000000001
$ws.Range("A10").Value = -21.604
$ws.Range("A12").Value = -21.606
$ws.Range("B15").Value = 4.867
$ws.Range("A18").Value = -22.192
$ws.Range("B20").Value = 7.031999999999999
$ws.Range("B29").Value = 5.645
$ws.Range("B30").Value = 6.02
$ws.Range("B31").Value = 6.367000000000001
$ws.Range("A37").Value = -19.92
$ws.Range("B40").Value = 9.327999999999999
$ws.Range("A55").Value = -21.795
$ws.Range("A68").Value = -21.507
$ws.Range("B68").Value = 5.881
$ws.Range("B76").Value = 6.308
$ws.Range("A77").Value = -20.637
$ws.Range("A78").Value = -20.126
$ws.Range("B87").Value = 4.836
$ws.Range("B88").Value = 5.058000000000001
$ws.Range("B96").Value = 6.692
$ws.Range("B98").Value = 5.095000000000001
$ws.Range("B101").Value = 7.782000000000001
$ws.Range("B102").Value = 7.747
